# Weekly refresh: insert the newest "Ciboulette" price record at the top of
# the data block (row 11), pushing all existing rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 11 (row 1 = header, rows 2-10 untouched,
# rows 11-111 shift down to 12-112).
$ws.Rows.Item(11).Insert()

# Populate the new row 11 with this week's record.
$ws.Cells.Item(11, 1).Value = 4
$ws.Cells.Item(11, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(11, 3).Value = "Los Lagos"
$ws.Cells.Item(11, 4).Value = 44462
$ws.Cells.Item(11, 5).Value = 10
$ws.Cells.Item(11, 6).Value = 100112039
$ws.Cells.Item(11, 7).Value = "Ciboulette"
$ws.Cells.Item(11, 8).Value = "Sin especificar"
$ws.Cells.Item(11, 9).Value = "Primera"
$ws.Cells.Item(11, 10).Value = 120
$ws.Cells.Item(11, 11).Value = 4000
$ws.Cells.Item(11, 12).Value = 4000
$ws.Cells.Item(11, 13).Value = 4000
$ws.Cells.Item(11, 14).Value = "`$/docena de atados"
$ws.Cells.Item(11, 15).Value = "Región Metropolitana"
$ws.Cells.Item(11, 16).Value = 1333
$ws.Cells.Item(11, 17).Value = 3
$ws.Cells.Item(11, 18).Value = "Hortaliza"
